{"js": "// Remove the standalone \"\u25b2\" run that precedes \"\u8868 8-2-10 \u5be9\u6838\u4eba\u54e1\"\n// in the last paragraph of the document body.\nconst searchResults = context.document.body.search(\"\u25b2\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the standalone \"\u25b2\" run that precedes \"\u8868 8-2-10 \u5be9\u6838\u4eba\u54e1\"\n# in the last paragraph of the document body.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"\u25b2\"\n$range.Find.Forward = $true\n$range.Find.Wrap = 1\n\nwhile ($range.Find.Execute()) {\n    $range.Delete()\n}\n"}
